# Error Calculations and Plots
# - Remove two rows that were dropped from the source data ("RM 232" and
#   "SC 92"), which shifts all subsequent rows up.
# - Toggle a handful of individual cell values between blank and a filled
#   numeric value (re-running the missing-data removal/imputation pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (originally row 26) and, after the shift, the
# "SC 92" row (originally row 28, now row 27).
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply the individual cell value changes (post row-shift row numbers).
$ws.Range("D2").Value = -13.5
$ws.Range("F3").ClearContents()
$ws.Range("F4").Value = 17.97
$ws.Range("F5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("D12").Value = -14.1
$ws.Range("D14").ClearContents()
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("F19").ClearContents()
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("F22").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("F23").Value = 16.48
$ws.Range("D24").ClearContents()
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("B30").Value = -19.7
$ws.Range("D31").Value = -13.7
$ws.Range("B32").ClearContents()
$ws.Range("D33").Value = -14.1
